$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3
$ws.Range("A2").Value = 11111111
$ws.Range("B2").Value = 123
$ws.Range("C2").Value = 5

$ws.Range("A3").Value = 11111111
$ws.Range("B3").Value = 125
$ws.Range("C3").Value = 1

# Add new row 4
$ws.Range("A4").Value = 22222222
$ws.Range("B4").Value = 125
$ws.Range("C4").Value = 1
